$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'54.250.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -8.60%  "

$ws.Range("D3").Value = "'2.404.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -15.33%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'463.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -8.51%  "

$ws.Range("D6").Value = "'130.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.28%  "

$ws.Range("D7").Value = "'0.994"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.61%  "

$ws.Range("D8").Value = "'0.486"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -7.94%  "

$ws.Range("D9").Value = "'2.418.72"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -14.70%  "

$ws.Range("D10").Value = "'0.0940"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -9.42%  "

$ws.Range("E11").Value = "  -9.95%  "

$ws.Range("E12").Value = "  -9.25%  "

$ws.Range("E13").Value = "  -4.08%  "

$ws.Range("D14").Value = "'2.815.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -15.56%  "

$ws.Range("D15").Value = "'54.016.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -9.05%  "

$ws.Range("E16").Value = "  -10.35%  "

$ws.Range("E17").Value = "  -5.50%  "

$ws.Range("D18").Value = "'2.413.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -15.01%  "

$ws.Range("E19").Value = "  -11.85%  "

$ws.Range("D20").Value = "'310.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -12.05%  "

$ws.Range("D21").Value = "'9.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -16.05%  "

$ws.Range("D22").Value = "'0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("D23").Value = "'5.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.32%  "

$ws.Range("D24").Value = "'5.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -14.31%  "

$ws.Range("D25").Value = "'56.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -10.66%  "

$ws.Range("E27").Value = "  -10.09%  "

$ws.Range("D28").Value = "'0.380"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -11.60%  "

$ws.Range("D29").Value = "'2.495.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -16.06%  "

$ws.Range("E30").Value = "  -4.87%  "

$ws.Range("D31").Value = "'0.996"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.34%  "

$ws.Range("D32").Value = "'0.0₃0705"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -13.96%  "

$ws.Range("D33").Value = "'149.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.73%  "

$ws.Range("D34").Value = "'17.59"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.54%  "

$ws.Range("E35").Value = "  -13.36%  "

$ws.Range("D36").Value = "'5.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.57%  "

$ws.Range("D37").Value = "'3.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -16.70%  "

$ws.Range("E38").Value = "  -8.84%  "

$ws.Range("D39").Value = "'0.797"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -14.85%  "

$ws.Range("D40").Value = "'33.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.06%  "

$ws.Range("D41").Value = "'0.992"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.60%  "

$ws.Range("D42").Value = "'0.598"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.15%  "

$ws.Range("D43").Value = "'0.0529"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.83%  "

$ws.Range("D44").Value = "'3.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.07%  "

$ws.Range("E45").Value = "  -1.94%  "

$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "'1.970.53"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -11.54%  "

$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "'1.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -10.50%  "

$ws.Range("D48").Value = "'0.0218"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.14%  "

$ws.Range("D49").Value = "'0.0864"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.77%  "

$ws.Range("D50").Value = "'4.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.53%  "

$ws.Range("D51").Value = "'16.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -15.97%  "

Write-Host "All updates applied."
